{"js": "// Replace the date in the title and the division problems in the table\n// with the values from the updated worksheet.\nconst replacements = [\n  [\"2024-10-07 Monday\", \"2024-10-08 Tuesday\"],\n  [\"29\u00f72=\", \"88\u00f79=\"],\n  [\"16\u00f75=\", \"29\u00f79=\"],\n  [\"38\u00f79=\", \"96\u00f72=\"],\n  [\"11\u00f77=\", \"82\u00f72=\"],\n  [\"86\u00f77=\", \"91\u00f77=\"],\n  [\"53\u00f74=\", \"21\u00f76=\"],\n  [\"74\u00f74=\", \"83\u00f74=\"],\n  [\"74\u00f73=\", \"76\u00f77=\"],\n  [\"18\u00f73=\", \"99\u00f79=\"],\n  [\"28\u00f76=\", \"53\u00f75=\"],\n  [\"35\u00f74=\", \"18\u00f77=\"],\n  [\"26\u00f77=\", \"25\u00f74=\"],\n  [\"44\u00f79=\", \"25\u00f76=\"],\n  [\"28\u00f73=\", \"37\u00f74=\"],\n  [\"46\u00f72=\", \"55\u00f73=\"],\n  [\"80\u00f72=\", \"17\u00f76=\"],\n  [\"83\u00f78=\", \"26\u00f79=\"],\n  [\"94\u00f75=\", \"40\u00f73=\"],\n  [\"26\u00f78=\", \"78\u00f72=\"],\n  [\"64\u00f73=\", \"78\u00f73=\"],\n  [\"93\u00f74=\", \"37\u00f75=\"],\n  [\"31\u00f73=\", \"98\u00f77=\"],\n  [\"97\u00f73=\", \"59\u00f72=\"],\n  [\"54\u00f74=\", \"78\u00f79=\"],\n  [\"41\u00f78=\", \"20\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and the division problems with the new values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2024-10-07 Monday\", \"2024-10-08 Tuesday\"),\n    @(\"29\u00f72=\", \"88\u00f79=\"),\n    @(\"16\u00f75=\", \"29\u00f79=\"),\n    @(\"38\u00f79=\", \"96\u00f72=\"),\n    @(\"11\u00f77=\", \"82\u00f72=\"),\n    @(\"86\u00f77=\", \"91\u00f77=\"),\n    @(\"53\u00f74=\", \"21\u00f76=\"),\n    @(\"74\u00f74=\", \"83\u00f74=\"),\n    @(\"74\u00f73=\", \"76\u00f77=\"),\n    @(\"18\u00f73=\", \"99\u00f79=\"),\n    @(\"28\u00f76=\", \"53\u00f75=\"),\n    @(\"35\u00f74=\", \"18\u00f77=\"),\n    @(\"26\u00f77=\", \"25\u00f74=\"),\n    @(\"44\u00f79=\", \"25\u00f76=\"),\n    @(\"28\u00f73=\", \"37\u00f74=\"),\n    @(\"46\u00f72=\", \"55\u00f73=\"),\n    @(\"80\u00f72=\", \"17\u00f76=\"),\n    @(\"83\u00f78=\", \"26\u00f79=\"),\n    @(\"94\u00f75=\", \"40\u00f73=\"),\n    @(\"26\u00f78=\", \"78\u00f72=\"),\n    @(\"64\u00f73=\", \"78\u00f73=\"),\n    @(\"93\u00f74=\", \"37\u00f75=\"),\n    @(\"31\u00f73=\", \"98\u00f77=\"),\n    @(\"97\u00f73=\", \"59\u00f72=\"),\n    @(\"54\u00f74=\", \"78\u00f79=\"),\n    @(\"41\u00f78=\", \"20\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
